$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Price" (D) and "Volume(1h)" (E) columns with refreshed crypto
# quotes. Values in column D that look like plain numbers are written with a
# leading apostrophe so Excel keeps them as literal text (preserving things
# like trailing zeros, e.g. "6.10") instead of silently converting them to
# numeric cell values.

$ws.Range('D2').Value = '58.008.20'
$ws.Range('E2').Value = '  -1.60%  '
$ws.Range('D3').Value = '3.098.53'
$ws.Range('E3').Value = '  -0.13%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').Value = '''527.07'
$ws.Range('E5').Value = '  +1.00%  '
$ws.Range('D6').Value = '''141.13'
$ws.Range('E6').Value = '  -1.98%  '
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('D8').Value = '3.098.74'
$ws.Range('E8').Value = '  -0.09%  '
$ws.Range('D9').Value = '''0.443'
$ws.Range('E9').Value = '  +1.00%  '
$ws.Range('D10').Value = '''7.17'
$ws.Range('E10').Value = '  -2.78%  '
$ws.Range('E11').Value = '  -0.93%  '
$ws.Range('E12').Value = '  +2.62%  '
$ws.Range('D13').Value = '3.630.04'
$ws.Range('E13').Value = '  -0.21%  '
$ws.Range('E14').Value = '  +2.53%  '
$ws.Range('D15').Value = '''25.74'
$ws.Range('E15').Value = '  -4.97%  '
$ws.Range('E16').Value = '  -1.24%  '
$ws.Range('D17').Value = '58.049.76'
$ws.Range('E17').Value = '  -1.54%  '
$ws.Range('D18').Value = '3.098.49'
$ws.Range('E18').Value = '  -0.30%  '
$ws.Range('D19').Value = '''6.10'
$ws.Range('E19').Value = '  -1.96%  '
$ws.Range('D20').Value = '''12.72'
$ws.Range('E20').Value = '  -2.45%  '
$ws.Range('D21').Value = '''7.98'
$ws.Range('E21').Value = '  -3.16%  '
$ws.Range('D22').Value = '''343.38'
$ws.Range('E22').Value = '  +0.16%  '
$ws.Range('D23').Value = '''0.999'
$ws.Range('E23').Value = '  -0.22%  '
$ws.Range('D24').Value = '''0.512'
$ws.Range('E24').Value = '  +0.64%  '
$ws.Range('E25').Value = '  +2.47%  '
$ws.Range('E26').Value = '  -1.11%  '
$ws.Range('E27').Value = '  +0.06%  '
$ws.Range('D28').Value = '0.0₃0921'
$ws.Range('E28').Value = '  -1.23%  '
$ws.Range('D29').Value = '''0.999'
$ws.Range('E29').Value = '  +0.00%  '
$ws.Range('D30').Value = '''6.40'
$ws.Range('E30').Value = '  -4.96%  '
$ws.Range('D31').Value = '''7.29'
$ws.Range('E31').Value = '  +0.24%  '
$ws.Range('E32').Value = '  +1.71%  '
$ws.Range('D33').Value = '''20.98'
$ws.Range('E33').Value = '  -0.30%  '
$ws.Range('D34').Value = '''1.19'
$ws.Range('E34').Value = '  -2.78%  '
$ws.Range('D35').Value = '''158.50'
$ws.Range('E35').Value = '  +2.09%  '
$ws.Range('D36').Value = '''4.64'
$ws.Range('E36').Value = '  -0.62%  '
$ws.Range('D37').Value = '''6.17'
$ws.Range('E37').Value = '  -0.39%  '
$ws.Range('D38').Value = '''26.17'
$ws.Range('E38').Value = '  -3.05%  '
$ws.Range('E39').Value = '  -4.98%  '
$ws.Range('D40').Value = '''0.0669'
$ws.Range('E40').Value = '  -2.51%  '
$ws.Range('D41').Value = '''4.04'
$ws.Range('E41').Value = '  +2.06%  '
$ws.Range('E42').Value = '  +6.16%  '
$ws.Range('D43').Value = '''0.685'
$ws.Range('E43').Value = '  +2.76%  '
$ws.Range('D44').Value = '3.138.00'
$ws.Range('E44').Value = '  -0.27%  '
$ws.Range('D45').Value = '''36.94'
$ws.Range('E45').Value = '  +0.19%  '
$ws.Range('E46').Value = '  -0.04%  '
$ws.Range('E47').Value = '  +2.03%  '
$ws.Range('D48').Value = '2.276.31'
$ws.Range('E48').Value = '  -0.37%  '
$ws.Range('D49').Value = '''0.991'
$ws.Range('E49').Value = '  +3.19%  '
$ws.Range('D50').Value = '''6.11'
$ws.Range('E50').Value = '  +1.47%  '
$ws.Range('D51').Value = '''20.51'
$ws.Range('E51').Value = '  -2.14%  '
